$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.446.34"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "3.455.28"
$ws.Range("E3").Value = "  +2.80%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'584.75"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").Value = "'179.78"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "'0.200"
$ws.Range("E9").Value = "  +7.74%  "

$ws.Range("D10").Value = "'0.589"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("D11").Value = "'49.36"
$ws.Range("E11").Value = "  +2.45%  "

$ws.Range("E12").Value = "  +2.73%  "

$ws.Range("D13").Value = "'683.59"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "4.007.88"
$ws.Range("E14").Value = "  +2.27%  "

$ws.Range("E15").Value = "  +2.41%  "

$ws.Range("D16").Value = "69.505.26"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").Value = "3.448.60"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").Value = "'17.93"
$ws.Range("E19").Value = "  +2.02%  "

$ws.Range("D20").Value = "'11.40"
$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("D22").Value = "'5.39"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").Value = "'17.12"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").Value = "'101.19"
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").Value = "'3.93"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").Value = "'2.72"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").Value = "'9.77"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("D28").Value = "'33.86"

$ws.Range("D29").Value = "'8.83"
$ws.Range("E29").Value = "  +2.75%  "

$ws.Range("D30").Value = "'6.92"
$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("D31").Value = "'3.77"
$ws.Range("E31").Value = "  +5.90%  "

$ws.Range("D32").Value = "'567.46"
$ws.Range("E32").Value = "  +2.43%  "

$ws.Range("D33").Value = "'11.09"
$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("D35").Value = "'58.11"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "3.643.97"
$ws.Range("E37").Value = "  -1.77%  "

$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "'35.33"
$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("E40").Value = "  +10.22%  "

$ws.Range("D41").Value = "'3.31"
$ws.Range("E41").Value = "  +3.35%  "

$ws.Range("D42").Value = "'2.71"
$ws.Range("E42").Value = "  +2.77%  "

$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.338"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0425"
$ws.Range("E45").Value = "  +2.42%  "

$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("E48").Value = "  +5.22%  "

$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").Value = "'131.43"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'2.69"
$ws.Range("E51").Value = "  +2.03%  "
